# AFDP-3899 Contacts Management - Organization - Organization ACL
#
# Adds a new "Organization" row to the Assignment Rules table on Sheet1,
# mirroring the existing "Complaint" / "Case File" rows (row 21 / row 20):
#   B22: Rule Name
#   C22: Type of Object to be Assigned
#   D22: Expression 1 (condition)
#   E22: Expression 2 (condition)
#   F22: Error message (action)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Clone the formatting (cell styles + row height) of the row above (the
# "Complaint" rule, row 21) onto the new row 22, the same way it was done
# for every previous rule row.
$ws.Range("B21:F21").Copy()
$ws.Range("B22:F22").PasteSpecial(-4122)  # xlPasteFormats
$ws.Rows.Item(22).RowHeight = 45

# Fill in the new rule's content.
$ws.Range("B22").Value = "Organization - Check participants list for NoAccess & Owner"
$ws.Range("C22").Value = "ORGANIZATION"
$ws.Range("D22").Value = "participants != null && participants.containsKey('No Access') && participants.containsKey('owner')"
$ws.Range("E22").Value = "participants['No Access'].contains(participants['owner'][0])"
$ws.Range("F22").Value = "Owners cannot be on the no-access list."

# Move the selection to the newly added cell, matching the saved view state.
$ws.Range("F22").Select()
